$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Summary")
$ws.Range("B3").Value = 0.01
$ws.Range("B6").Value = 107190.1972019555
$ws.Range("B7").Value = 12382564.12454636
$ws.Range("B8").Value = 20513624.20839428
$ws.Range("B10").Value = 4288679.125387715

$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Range("B2").Value = 491565.8032302229
$ws.Range("C2").Value = 491565.8032302228
$ws.Range("E2").Value = 79335.53193155181
$ws.Range("F2").Value = 67915.6684446861
$ws.Range("H2").Value = 106263.3325344874
$ws.Range("I2").Value = 85771.85995002884
$ws.Range("J2").Value = 86568.55741505309
$ws.Range("K2").Value = 108380.894185064
$ws.Range("L2").Value = 97614.64873949091
$ws.Range("M2").Value = 65003.6223198295
$ws.Range("N2").Value = 79907.52366022246
$ws.Range("O2").Value = 89394.47967254139
$ws.Range("P2").Value = 64519.26044478228
$ws.Range("E4").Value = 19377.29391253668
$ws.Range("F4").Value = 7481.602780384896
$ws.Range("H4").Value = 47427.08620726129
$ws.Range("I4").Value = 26081.80226511689
$ws.Range("J4").Value = 26911.6954578505
$ws.Range("K4").Value = 49632.87959327855
$ws.Range("L4").Value = 38418.04058747327
$ws.Range("M4").Value = 4448.221400325934
$ws.Range("N4").Value = 19973.1186299019
$ws.Range("O4").Value = 29855.3644760675
$ws.Range("P4").Value = 3943.677780485084
$ws.Range("B6").Value = -9721.89248018945
$ws.Range("C6").Value = -9721.892480189508
$ws.Range("D6").Value = -9721.892480189217
$ws.Range("E6").Value = -93973.77116878566
$ws.Range("F6").Value = 39126.22883121439
$ws.Range("G6").Value = 39126.22883121438
$ws.Range("H6").Value = 39126.22883121436
$ws.Range("I6").Value = 39126.22883121441
$ws.Range("J6").Value = 39126.22883121439
$ws.Range("K6").Value = 39126.22883121436
$ws.Range("L6").Value = 39126.22883121435
$ws.Range("M6").Value = 39126.22883121439
$ws.Range("N6").Value = 39126.22883121442
$ws.Range("O6").Value = 39126.22883121438
$ws.Range("P6").Value = 39126.22883121438

$ws = $wb.Worksheets.Item("DG Dispatch")
$ws.Range("J11").Value = 124.5190384721106
$ws.Range("K11").Value = 135.370731907559
$ws.Range("I12").Value = 10.12574714858493
$ws.Range("J12").Value = 93.17061249236157
$ws.Range("Q12").Value = 94.49434172313325
$ws.Range("R12").Value = 45.52166981132082
$ws.Range("Q14").Value = 150.3839754851235
$ws.Range("J20").Value = 124.5190384721106
$ws.Range("K20").Value = 135.370731907559
$ws.Range("M20").Value = 113.4004983079896
$ws.Range("N20").Value = 110.5750244233121
$ws.Range("O20").Value = 117.8828208804077
$ws.Range("P20").Value = 135.4597561231036
$ws.Range("Q20").Value = 150.3839754851235
$ws.Range("R20").Value = 65.71641987298243
$ws.Range("J23").Value = 124.5190384721106
$ws.Range("K23").Value = 135.370731907559
$ws.Range("L23").Value = 130.6648563030561
$ws.Range("I24").Value = 10.12574714858493
$ws.Range("R24").Value = 45.52166981132082
$ws.Range("N25").Value = 81.96869489115805
$ws.Range("O25").Value = 96.22962838366004
$ws.Range("P25").Value = 101.5955875616828
$ws.Range("Q25").Value = 65.34295837775146
$ws.Range("J26").Value = 124.5190384721106
$ws.Range("K26").Value = 135.370731907559
$ws.Range("L26").Value = 130.6648563030561
$ws.Range("Q26").Value = 150.3839754851235
$ws.Range("J29").Value = 124.5190384721106
$ws.Range("K29").Value = 135.370731907559
$ws.Range("L29").Value = 130.6648563030561
$ws.Range("M29").Value = 113.4004983079896
$ws.Range("O29").Value = 117.8828208804077
$ws.Range("P29").Value = 135.4597561231036
$ws.Range("Q29").Value = 150.3839754851235
$ws.Range("R29").Value = 65.71641987298243
$ws.Range("R30").Value = 45.52166981132082
$ws.Range("J32").Value = 124.5190384721106
$ws.Range("K32").Value = 135.370731907559
$ws.Range("L32").Value = 130.6648563030561
$ws.Range("M32").Value = 113.4004983079896
$ws.Range("O32").Value = 117.8828208804077
$ws.Range("Q32").Value = 150.3839754851235
$ws.Range("I36").Value = 10.12574714858493
$ws.Range("J36").Value = 93.17061249236157
$ws.Range("R36").Value = 45.52166981132082
$ws.Range("J37").Value = 33.63624132272333
$ws.Range("J38").Value = 124.5190384721106
$ws.Range("I39").Value = 10.12574714858493
$ws.Range("J39").Value = 93.17061249236157
$ws.Range("R39").Value = 45.52166981132082
$ws.Range("J40").Value = 33.63624132272333
$ws.Range("K40").Value = 94.30397654773019
$ws.Range("L40").Value = 90.4687457914608
$ws.Range("M40").Value = 92.09541281912071
$ws.Range("N40").Value = 81.96869489115805
$ws.Range("O40").Value = 96.22962838366004
$ws.Range("P40").Value = 101.5955875616828
$ws.Range("Q40").Value = 65.34295837775146
$ws.Range("L41").Value = 130.6648563030561
$ws.Range("O41").Value = 117.8828208804077
$ws.Range("P41").Value = 135.4597561231036
$ws.Range("Q41").Value = 150.3839754851235
$ws.Range("R41").Value = 65.71641987298243
$ws.Range("I45").Value = 10.12574714858493
$ws.Range("J45").Value = 93.17061249236157
$ws.Range("R45").Value = 45.52166981132082

$ws = $wb.Worksheets.Item("Unmet Demand")
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("I12").Value = 77.12765456497084
$ws.Range("J12").Value = 0
$ws.Range("Q12").Value = 0
$ws.Range("R12").Value = 78.03303713061706
$ws.Range("Q14").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = 0
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = 0
$ws.Range("Q20").Value = 0
$ws.Range("R20").Value = 108.0327934026353
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("I24").Value = 77.12765456497084
$ws.Range("R24").Value = 78.03303713061706
$ws.Range("N25").Value = 0
$ws.Range("O25").Value = 0
$ws.Range("P25").Value = 0
$ws.Range("Q25").Value = 61.14583096471014
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("Q26").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 0
$ws.Range("O29").Value = 0
$ws.Range("P29").Value = 0
$ws.Range("Q29").Value = 0
$ws.Range("R29").Value = 108.0327934026353
$ws.Range("R30").Value = 78.03303713061706
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = 0
$ws.Range("O32").Value = 0
$ws.Range("Q32").Value = 0
$ws.Range("I36").Value = 77.12765456497084
$ws.Range("J36").Value = 0
$ws.Range("R36").Value = 78.03303713061706
$ws.Range("J37").Value = 72.23757736389061
$ws.Range("J38").Value = 0
$ws.Range("I39").Value = 77.12765456497084
$ws.Range("J39").Value = 0
$ws.Range("R39").Value = 78.03303713061706
$ws.Range("J40").Value = 72.23757736389061
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = 0
$ws.Range("N40").Value = 0
$ws.Range("O40").Value = 0
$ws.Range("P40").Value = 0
$ws.Range("Q40").Value = 61.14583096471014
$ws.Range("L41").Value = 0
$ws.Range("O41").Value = 0
$ws.Range("P41").Value = 0
$ws.Range("Q41").Value = 0
$ws.Range("R41").Value = 108.0327934026353
$ws.Range("I45").Value = 77.12765456497084
$ws.Range("J45").Value = 0
$ws.Range("R45").Value = 0

$ws = $wb.Worksheets.Item("Household Surplus")
$ws.Range("B5").Value = 171437.0634480109
$ws.Range("B6").Value = 137653.3006326999
$ws.Range("B8").Value = 251098.4735650287
$ws.Range("B9").Value = 190477.8671693386
$ws.Range("B10").Value = 192834.7638367022
$ws.Range("B11").Value = 257362.9267813178
$ws.Range("B12").Value = 225512.7840048308
$ws.Range("B13").Value = 129038.4975133324
$ws.Range("B14").Value = 173129.2056453282
$ws.Range("B15").Value = 201194.7838484384
$ws.Range("B16").Value = 127605.5936329844
